# Updated cryptos list on Fri Apr 19 14:43:02 UTC 2024 with GitHub Actions
#
# Refreshes the coinranking.com snapshot captured in the worksheet: new
# "Price" (col D) and "Volume(1h)" (col E) readings for each coin, plus a
# rank swap between dogwifhat and VeChain (rows 39-40) now that dogwifhat
# overtook VeChain.
#
# Price/percentage cells are stored as literal text in the source feed
# (e.g. "64.243.04", "  +0.66%  "), so each write temporarily forces
# NumberFormat "@" (text) to stop Excel's automatic number coercion, then
# restores the cell's original Style afterwards to avoid leaving any
# formatting residue behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $cell = $ws.Range($range)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue 'D2' '64.243.04'
Set-TextValue 'E2' '  +0.66%  '
Set-TextValue 'D3' '3.072.76'
Set-TextValue 'E3' '  -0.36%  '
Set-TextValue 'D5' '556.74'
Set-TextValue 'E5' '  +0.68%  '
Set-TextValue 'D6' '143.67'
Set-TextValue 'E6' '  +3.49%  '
Set-TextValue 'E7' '  +0.15%  '
Set-TextValue 'D8' '3.069.75'
Set-TextValue 'E8' '  -0.32%  '
Set-TextValue 'D9' '0.496'
Set-TextValue 'E9' '  -0.75%  '
Set-TextValue 'D10' '6.32'
Set-TextValue 'E10' '  +2.23%  '
Set-TextValue 'E11' '  -0.02%  '
Set-TextValue 'D12' '0.469'
Set-TextValue 'E12' '  +2.92%  '
Set-TextValue 'D13' '0.0000228'
Set-TextValue 'E13' '  +0.21%  '
Set-TextValue 'D14' '35.12'
Set-TextValue 'E14' '  +0.40%  '
Set-TextValue 'D15' '3.590.00'
Set-TextValue 'E15' '  +0.51%  '
Set-TextValue 'D16' '64.297.05'
Set-TextValue 'E16' '  +0.99%  '
Set-TextValue 'D17' '3.079.91'
Set-TextValue 'E17' '  -0.13%  '
Set-TextValue 'D18' '0.110'
Set-TextValue 'E18' '  +0.47%  '
Set-TextValue 'D19' '6.75'
Set-TextValue 'E19' '  -0.01%  '
Set-TextValue 'D20' '474.10'
Set-TextValue 'E20' '  -2.64%  '
Set-TextValue 'D21' '13.86'
Set-TextValue 'E21' '  +2.36%  '
Set-TextValue 'D22' '0.680'
Set-TextValue 'E22' '  -0.82%  '
Set-TextValue 'E23' '  +4.24%  '
Set-TextValue 'D24' '13.50'
Set-TextValue 'E24' '  +8.28%  '
Set-TextValue 'D25' '80.78'
Set-TextValue 'E25' '  -0.90%  '
Set-TextValue 'D26' '1.00'
Set-TextValue 'E26' '  +0.12%  '
Set-TextValue 'D27' '2.77'
Set-TextValue 'E27' '  +0.42%  '
Set-TextValue 'D28' '8.12'
Set-TextValue 'E28' '  +1.04%  '
Set-TextValue 'E29' '  +2.49%  '
Set-TextValue 'D30' '0.998'
Set-TextValue 'E30' '  -0.40%  '
Set-TextValue 'D31' '25.93'
Set-TextValue 'E31' '  -0.55%  '
Set-TextValue 'E32' '  +0.25%  '
Set-TextValue 'D33' '2.48'
Set-TextValue 'E33' '  +2.23%  '
Set-TextValue 'D34' '5.57'
Set-TextValue 'E34' '  -4.39%  '
Set-TextValue 'E35' '  +1.83%  '
Set-TextValue 'D36' '54.58'
Set-TextValue 'E36' '  -2.02%  '
Set-TextValue 'D37' '465.14'
Set-TextValue 'E37' '  -1.56%  '
Set-TextValue 'D38' '0.0830'
Set-TextValue 'E38' '  +1.33%  '
Set-TextValue 'B39' 'dogwifhat'
Set-TextValue 'C39' 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue 'D39' '2.96'
Set-TextValue 'E39' '  +14.64%  '
Set-TextValue 'B40' 'VeChain'
Set-TextValue 'C40' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D40' '0.0405'
Set-TextValue 'E40' '  +1.96%  '
Set-TextValue 'D41' '2.977.12'
Set-TextValue 'E41' '  -6.79%  '
Set-TextValue 'D42' '8.23'
Set-TextValue 'E42' '  +0.04%  '
Set-TextValue 'E43' '  -5.66%  '
Set-TextValue 'D44' '28.46'
Set-TextValue 'E44' '  +1.91%  '
Set-TextValue 'E45' '  +1.88%  '
Set-TextValue 'E47' '  +3.73%  '
Set-TextValue 'E48' '  +1.38%  '
Set-TextValue 'D49' '0.0₃0519'
Set-TextValue 'E49' '  +0.69%  '
Set-TextValue 'D50' '117.61'
Set-TextValue 'E50' '  +1.09%  '
Set-TextValue 'D51' '2.06'
Set-TextValue 'E51' '  -0.34%  '
